$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a brand-new column before column A ("TestCases" data provider column).
# This shifts the existing Tasks/Completed columns from A:B to B:C, carrying
# along their column widths, cell styles and values automatically.
$ws.Columns("A:A").Insert()

# Header for the new column, bold like the other headers.
$ws.Range("A1").Value = "TestCases"
$ws.Range("A1").Font.Bold = $true

# TestCase id values for each data row.
$ws.Range("A2").Value = "TC_01"
$ws.Range("A3").Value = "TC_02"
$ws.Range("A4").Value = "TC_02"
$ws.Range("A5").Value = "TC_02"
$ws.Range("A6").Value = "TC_02"
$ws.Range("A7").Value = "TC_02"
$ws.Range("A8").Value = "TC_03"
$ws.Range("A9").Value = "TC_03"
$ws.Range("A10").Value = "TC_03"
$ws.Range("A11").Value = "TC_04"
$ws.Range("A12").Value = "TC_04"
$ws.Range("A13").Value = "TC_05"
$ws.Range("A14").Value = "TC_06"
$ws.Range("A15").Value = "TC_06"

# Carry the formatted style of row 8 (B8:C8) down through the rest of the
# table (rows 9-15) so every data cell -- including the brand new rows -- has
# the same look as the existing rows, then fill in the actual values.
$ws.Range("B8:C8").Copy()
$ws.Range("B9:C15").PasteSpecial(-4122)

$ws.Range("B9").Value = "Component Testing"
$ws.Range("B10").Value = "Component Testing"

$ws.Range("B11").Value = ""
$ws.Range("C11").Value = "Requirement Gathering"

$ws.Range("B12").Value = ""
$ws.Range("C12").Value = "Requirement Analysis"

$ws.Range("C13").Value = ""
$ws.Range("B13").Value = "Integration Testing"

$ws.Range("C14").Value = ""
$ws.Range("B14").Value = "System Testing"

$ws.Range("C15").Value = ""
$ws.Range("B15").Value = "Production Support"

# Match the final selection left behind in the saved workbook.
$ws.Range("A14:A15").Select() | Out-Null
